$d = $word.ActiveDocument

# Locate the paragraph that contains "LOQ4086" (the requirement line that stays),
# and the paragraph that contains the copyright/footer text "Creative Commons
# Attribution" (the last paragraph to be removed). Everything between them
# (an empty paragraph, the "Ver no Jupiter..." paragraph, and the "(c) 2020..."
# paragraph itself) is deleted, per the commit's removal of the site-footer
# boilerplate that used to be appended after each course requirement listing.

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOQ4086*") {
        $startIdx = $i + 1
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1 -and $endIdx -ge $startIdx) {
    $pStart = $d.Paragraphs.Item($startIdx)
    $pEnd = $d.Paragraphs.Item($endIdx)
    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}
